$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing trade data for row 2 (bollinger bands calc fix)
$ws.Range("B2").Value = $true
$ws.Range("E2").Value = 108.91
$ws.Range("F2").Value = 0.45194613539936812
$ws.Range("G2").Value = $false

# Add new row 3 with updated principle value
$ws.Range("C3").Value = 10045.19

# Recompute best-fit column width for column C given the new content
# (stored OOXML "width" = ColumnWidth + 0.8333333333333334, so this yields
# the best-fit width of 9 that Excel computes for "10045.19")
$ws.Columns.Item(3).ColumnWidth = 8.166666666666666
